$d = $word.ActiveDocument

# Unicode helper characters used in the new prose (curly quotes / apostrophe)
$LQ = [char]0x201C   # “
$RQ = [char]0x201D   # ”
$AP = [char]0x2019   # ’

# =====================================================================
# Change 1: Paragraph 1 — "Forest" becomes a bold 3-run title
# =====================================================================
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range

$p1seg1 = "-- You choose to do "
$p1seg2 = "enter the fantasy forest"
$p1seg3 = " --"

$r1.Text = $p1seg1 + $p1seg2 + $p1seg3

$p1Start = $d.Paragraphs.Item(1).Range.Start
$b1 = $p1Start
$e1 = $b1 + $p1seg1.Length
$b2 = $e1
$e2 = $b2 + $p1seg2.Length
$b3 = $e2
$e3 = $b3 + $p1seg3.Length

$rg = $d.Range($b1, $e1)
$rg.Bold = 1
$rg.BoldBi = 1

$rg = $d.Range($b2, $e2)
$rg.Bold = 1
$rg.BoldBi = 1

$rg = $d.Range($b3, $e3)
$rg.Bold = 1
$rg.BoldBi = 1

# =====================================================================
# Change 2: Paragraph 3 — append the narrator giving starting items,
# then relocate the "_GoBack" bookmark to the end of this paragraph.
# =====================================================================
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$p3TextEnd = $r3.End - 1

$p3seg1 = " " + $LQ + "It" + $AP + "s been a while since someone chose to go in so quickly, here" + $AP + "s your starting items," + $RQ + " Narrator waves her hands and a short sword, a small bag, and a couple of potions materialize onto your hands."
$p3seg2 = " You peer into the bag to find a small amount of sparkly gold. `"Tha"
$p3seg3 = "t bag contains a total of 100 G," + $RQ + " she states."

$r3.InsertAfter($p3seg1 + $p3seg2 + $p3seg3)

$b1 = $p3TextEnd
$e1 = $b1 + $p3seg1.Length
$b2 = $e1
$e2 = $b2 + $p3seg2.Length

# Force the appended text to land in 3 distinct runs (matching the
# source) by toggling Bold on/off at each boundary — this is a no-op
# visually (default is already "not bold") but it breaks run-merging.
$rg = $d.Range($b1, $e2)
$rg.Bold = 1
$rg.Bold = 0

$rg = $d.Range($b2, $e2)
$rg.Bold = 1
$rg.Bold = 0

# Relocate the singleton "_GoBack" bookmark to the very end of
# paragraph 3. A collapsed range sitting exactly on the paragraph
# mark boundary is ambiguous, so append a throwaway character, anchor
# the bookmark just before it, then remove the throwaway character.
$p3b = $d.Paragraphs.Item(3)
$p3b.Range.InsertAfter("Z")
$p3c = $d.Paragraphs.Item(3)
$bmPos = $p3c.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$p3d = $d.Paragraphs.Item(3)
$delPos = $p3d.Range.End - 2
$d.Range($delPos, $delPos + 1).Delete()

# =====================================================================
# Change 3: Paragraph 9 — rewrite the ending into 6 runs
# =====================================================================
$p9 = $d.Paragraphs.Item(9)
$r9 = $p9.Range
$p9Start = $r9.Start

$p9seg1 = "As you fall"
$p9seg2 = " you hear, `"You were too overconfident a"
$p9seg3 = "nd look where that lead you,"
$p9seg4 = "`" before you fall into the mouth of the mythical hungry beast and "
$p9seg5 = "let out a bloodcurdling scream"
$p9seg6 = ". "

$r9.Text = $p9seg1 + $p9seg2 + $p9seg3 + $p9seg4 + $p9seg5 + $p9seg6

$b1 = $p9Start
$e1 = $b1 + $p9seg1.Length
$b2 = $e1
$e2 = $b2 + $p9seg2.Length
$b3 = $e2
$e3 = $b3 + $p9seg3.Length
$b4 = $e3
$e4 = $b4 + $p9seg4.Length
$b5 = $e4
$e5 = $b5 + $p9seg5.Length
$b6 = $e5
$e6 = $b6 + $p9seg6.Length

$rg = $d.Range($b2, $e6)
$rg.Bold = 1
$rg.Bold = 0

$rg = $d.Range($b3, $e6)
$rg.Bold = 1
$rg.Bold = 0

$rg = $d.Range($b4, $e6)
$rg.Bold = 1
$rg.Bold = 0

$rg = $d.Range($b5, $e6)
$rg.Bold = 1
$rg.Bold = 0

$rg = $d.Range($b6, $e6)
$rg.Bold = 1
$rg.Bold = 0

# =====================================================================
# Change 4: Paragraph 11 — merge the two runs left by the bookmark's
# old position back into a single run (bookmark itself already moved
# away in Change 2, which implicitly removed it from here).
# =====================================================================
$p11 = $d.Paragraphs.Item(11)
$p11FullText = $p11.Range.Text
$p11FullText = $p11FullText.Substring(0, $p11FullText.Length - 1)

$rngFull = $d.Range($p11.Range.Start, $p11.Range.End - 1)
$rngFull.Text = $p11FullText + "Z"
$p11b = $d.Paragraphs.Item(11)
$rngFull2 = $d.Range($p11b.Range.Start, $p11b.Range.End - 1)
$rngFull2.Text = $p11FullText
